$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the date value in A92 (R script re-ran with corrected timestamp) ---
$ws.Range("A92").Value = 45447.2916666667

# --- Append a new row 93 with the latest R-script results ---

# A93: date/time value, formatted like the other date cells in column A
$ws.Range("A93").Value = 45448.5009259259
$ws.Range("A92").Copy() | Out-Null
$ws.Range("A93").PasteSpecial(-4122) | Out-Null

$ws.Range("B93").Value = 900
$ws.Range("C93").Value = 2
$ws.Range("D93").Value = 2
$ws.Range("E93").Value = 2
$ws.Range("F93").Value = 2

# G93: adj_close stored as text "2" (matches the rest of column G), without
# leaving the cell tagged with a "Text" number format.
$ws.Range("G93").NumberFormat = "@"
$ws.Range("G93").Value = "2"
$ws.Range("B92").Copy() | Out-Null
$ws.Range("G93").PasteSpecial(-4122) | Out-Null

$ws.Range("H93").Value = "KK.MI"
